# Edit: break out stock.yaml completed
# 1) Append 12 new "day" rows (rows 555-566) to the "day" sheet, replicating
#    the intraday snapshot appended for the 16/09/2024 11:34:43 run.
# 2) Fix 31 mis-typed "bsecode" cells on the "week" sheet (rows 290, 292-321)
#    that had been written as text instead of numbers.

$wb = $excel.ActiveWorkbook
$wsDay = $wb.Worksheets.Item("day")
$wsWeek = $wb.Worksheets.Item("week")

$newRows = @(
  @{ A=1; B="OFSS"; C="Oracle Financial Services Software Limited"; D="532466"; E=-0.19; F=12239.1; G=377298; H="day"; I="16/09/2024 11:34:43" }
  @{ A=2; B="EICHERMOT"; C="Eicher Motors Limited"; D="505200"; E=0.84; F=4899.2; G=453853; H="day"; I="16/09/2024 11:34:43" }
  @{ A=3; B="TVSMOTOR"; C="Tvs Motor Company Limited"; D="532343"; E=0.46; F=2841.75; G=237258; H="day"; I="16/09/2024 11:34:43" }
  @{ A=4; B="INFY"; C="Infosys Limited"; D="500209"; E=0.32; F=1950.25; G=1570183; H="day"; I="16/09/2024 11:34:43" }
  @{ A=5; B="CIPLA"; C="Cipla Limited"; D="500087"; E=-0.02; F=1659.4; G=683445; H="day"; I="16/09/2024 11:34:43" }
  @{ A=6; B="TATACONSUM"; C="TATA Consumer Products Ltd"; D="500800"; E=0.68; F=1218.5; G=984662; H="day"; I="16/09/2024 11:34:43" }
  @{ A=7; B="GUJGASLTD"; C="Gujarat Gas Limited"; D="539336"; E=-0.43; F=630.3; G=642352; H="day"; I="16/09/2024 11:34:43" }
  @{ A=8; B="IGL"; C="Indraprastha Gas Limited"; D="532514"; E=2.46; F=529.85; G=3037121; H="day"; I="16/09/2024 11:34:43" }
  @{ A=9; B="CROMPTON"; C="Crompton Greaves Consumer Electricals Limited"; D="539876"; E=-0.85; F=447.9; G=2595847; H="day"; I="16/09/2024 11:34:43" }
  @{ A=10; B="VEDL"; C="Vedanta Limited"; D="500295"; E=-1.71; F=446.3; G=34590598; H="day"; I="16/09/2024 11:34:43" }
  @{ A=11; B="BPCL"; C="Bharat Petroleum Corporation Limited"; D="500547"; E=-0.48; F=340.65; G=4850244; H="day"; I="16/09/2024 11:34:43" }
  @{ A=12; B="IDFC"; C="Idfc Limited"; D="532659"; E=-0.37; F=110.99; G=3348343; H="day"; I="16/09/2024 11:34:43" }
)

$startRow = 555
$r = $startRow
foreach ($row in $newRows) {
    $wsDay.Cells.Item($r, 1).Value = $row.A
    $wsDay.Cells.Item($r, 2).Value = $row.B
    $wsDay.Cells.Item($r, 3).Value = $row.C

    # bsecode must stay textual (matches existing column formatting for
    # these rows), so force text before assigning the numeric-looking code.
    $dCell = $wsDay.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row.D

    $wsDay.Cells.Item($r, 5).Value = $row.E
    $wsDay.Cells.Item($r, 6).Value = $row.F
    $wsDay.Cells.Item($r, 7).Value = $row.G
    $wsDay.Cells.Item($r, 8).Value = $row.H
    $wsDay.Cells.Item($r, 9).Value = $row.I

    $r++
}

$weekFix = @(
  @{ Row=290; Value=500387 }
  @{ Row=292; Value=532466 }
  @{ Row=293; Value=540005 }
  @{ Row=294; Value=539448 }
  @{ Row=295; Value=505200 }
  @{ Row=296; Value=541154 }
  @{ Row=297; Value=500480 }
  @{ Row=298; Value=542726 }
  @{ Row=299; Value=533150 }
  @{ Row=300; Value=500325 }
  @{ Row=301; Value=500520 }
  @{ Row=302; Value=532187 }
  @{ Row=303; Value=532215 }
  @{ Row=304; Value=532321 }
  @{ Row=305; Value=500228 }
  @{ Row=306; Value=511196 }
  @{ Row=307; Value=532733 }
  @{ Row=308; Value=500440 }
  @{ Row=309; Value=532400 }
  @{ Row=310; Value=500670 }
  @{ Row=311; Value=512070 }
  @{ Row=312; Value=524208 }
  @{ Row=313; Value=507685 }
  @{ Row=314; Value=532810 }
  @{ Row=315; Value=500295 }
  @{ Row=316; Value=513599 }
  @{ Row=317; Value=532134 }
  @{ Row=318; Value=540691 }
  @{ Row=319; Value=531213 }
  @{ Row=320; Value=517334 }
  @{ Row=321; Value=532659 }
)

foreach ($fix in $weekFix) {
    # Was stored as text (e.g. "500387"); re-assigning a numeric literal
    # converts the cell to a real number, matching the other bsecode cells.
    $wsWeek.Cells.Item($fix.Row, 4).Value = $fix.Value
}
